$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3652
$ws.Range("H92").Value = 332.42307
$ws.Range("I92").Value = 322.9
$ws.Range("J92").Value = 364.16666
$ws.Range("K92").Value = 322.9
$ws.Range("L92").Value = 364.16666
$ws.Range("M92").Value = 925.1
$ws.Range("N92").Value = -2860.16666
$ws.Range("H97").Value = 371.5
$ws.Range("I97").Value = 352.5
$ws.Range("J97").Value = 375.3
$ws.Range("K97").Value = 1057.5
$ws.Range("L97").Value = 1125.9
$ws.Range("M97").Value = -561.5
$ws.Range("N97").Value = -2117.9
$ws.Range("H99").Value = 488.25
$ws.Range("I99").Value = 317
$ws.Range("J99").Value = 773.6667
$ws.Range("K99").Value = 951
$ws.Range("L99").Value = 2321.0001
$ws.Range("M99").Value = 547
$ws.Range("N99").Value = -5317.0001
$ws.Range("H100").Value = 1514.4
$ws.Range("J100").Value = 1785.8334
$ws.Range("L100").Value = 1785.8334
$ws.Range("N100").Value = -2867.8334
$ws.Range("H101").Value = 7501.1
$ws.Range("I101").Value = 528
$ws.Range("J101").Value = 13206.363
$ws.Range("K101").Value = 1584
$ws.Range("L101").Value = 39619.089
$ws.Range("M101").Value = 38
$ws.Range("N101").Value = -42863.089
$ws.Range("H106").Value = 3334.762
$ws.Range("I106").Value = 2279.111
$ws.Range("K106").Value = 2279.111
$ws.Range("M106").Value = -1648.111
$ws.Range("H132").Value = 2373.5217
$ws.Range("I132").Value = 2089.1296
$ws.Range("K132").Value = 6267.388800000001
$ws.Range("M132").Value = -3737.388800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1913.5714
$ws.Range("I74").Value = 1919.9032
$ws.Range("J74").Value = 1895.7273
$ws.Range("K74").Value = 1919.9032
$ws.Range("L74").Value = 1895.7273
$ws.Range("M74").Value = -1045.9032
$ws.Range("N74").Value = -3643.7273
$ws.Range("H77").Value = 1913.5714
$ws.Range("I77").Value = 1919.9032
$ws.Range("J77").Value = 1895.7273
$ws.Range("K77").Value = 9599.516
$ws.Range("L77").Value = 9478.636500000001
$ws.Range("M77").Value = -5231.516
$ws.Range("N77").Value = -18214.6365
$ws.Range("H97").Value = 1220
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 1850
$ws.Range("K97").Value = 800
$ws.Range("L97").Value = 1850
$ws.Range("M97").Value = -304
$ws.Range("N97").Value = -2842
$ws.Range("H102").Value = 1165.9
$ws.Range("I102").Value = 1028.0588
$ws.Range("J102").Value = 1947
$ws.Range("K102").Value = 1028.0588
$ws.Range("L102").Value = 1947
$ws.Range("M102").Value = 593.9412
$ws.Range("N102").Value = -5191
$ws.Range("H132").Value = 1307.9395
$ws.Range("I132").Value = 1095.341
$ws.Range("J132").Value = 1733.1364
$ws.Range("K132").Value = 3286.023
$ws.Range("L132").Value = 5199.4092
$ws.Range("M132").Value = -756.0229999999997
$ws.Range("N132").Value = -10259.4092

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 10065
$ws.Range("J49").Value = 10065
$ws.Range("L49").Value = 10065
$ws.Range("N49").Value = -10543
$ws.Range("H54").Value = 2972.6667
$ws.Range("I54").Value = 1567.2
$ws.Range("K54").Value = 1567.2
$ws.Range("M54").Value = -1083.2
$ws.Range("H64").Value = 431.8889
$ws.Range("J64").Value = 443.85715
$ws.Range("L64").Value = 443.85715
$ws.Range("N64").Value = -893.85715
$ws.Range("H67").Value = 431.8889
$ws.Range("J67").Value = 443.85715
$ws.Range("L67").Value = 443.85715
$ws.Range("N67").Value = -2003.85715
$ws.Range("H86").Value = 1988.3684
$ws.Range("I86").Value = 2059.923
$ws.Range("J86").Value = 1833.3334
$ws.Range("K86").Value = 2059.923
$ws.Range("L86").Value = 1833.3334
$ws.Range("M86").Value = -936.9229999999998
$ws.Range("N86").Value = -4079.3334
$ws.Range("H89").Value = 1988.3684
$ws.Range("I89").Value = 2059.923
$ws.Range("J89").Value = 1833.3334
$ws.Range("K89").Value = 10299.615
$ws.Range("L89").Value = 9166.666999999999
$ws.Range("M89").Value = -4683.614999999998
$ws.Range("N89").Value = -20398.667
$ws.Range("H94").Value = 487.47058
$ws.Range("I94").Value = 381.8
$ws.Range("J94").Value = 638.4286
$ws.Range("K94").Value = 381.8
$ws.Range("L94").Value = 638.4286
$ws.Range("M94").Value = 69.19999999999999
$ws.Range("N94").Value = -1540.4286
$ws.Range("H99").Value = 7504.5557
$ws.Range("I99").Value = 13218.625
$ws.Range("K99").Value = 13218.625
$ws.Range("M99").Value = -11720.625
$ws.Range("H134").Value = 1278.9855
$ws.Range("I134").Value = 1042.6786
$ws.Range("K134").Value = 3128.0358
$ws.Range("M134").Value = -593.0357999999997

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1453.5217
$ws.Range("I58").Value = 784.97437
$ws.Range("K58").Value = 784.97437
$ws.Range("M58").Value = -581.97437
$ws.Range("H122").Value = 2037.4546
$ws.Range("I122").Value = 1426.5
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 4279.5
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -1829.5
$ws.Range("N122").Value = -15900.0001
$ws.Range("H132").Value = 1483.6415
$ws.Range("I132").Value = 1086.7106
$ws.Range("J132").Value = 2489.2
$ws.Range("K132").Value = 3260.1318
$ws.Range("L132").Value = 7467.599999999999
$ws.Range("M132").Value = -730.1318000000001
$ws.Range("N132").Value = -12527.6
$ws.Range("H134").Value = 1957.2325
$ws.Range("I134").Value = 1325.5758
$ws.Range("J134").Value = 4041.7
$ws.Range("K134").Value = 3976.7274
$ws.Range("L134").Value = 12125.1
$ws.Range("M134").Value = -1441.7274
$ws.Range("N134").Value = -17195.1
$ws.Range("H136").Value = 1453.5217
$ws.Range("I136").Value = 784.97437
$ws.Range("K136").Value = 2354.92311
$ws.Range("M136").Value = 195.0768899999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4310789
$ws.Range("I113").Value = 17241730
$ws.Range("J113").Value = 475.33334
$ws.Range("K113").Value = 51725190
$ws.Range("L113").Value = 1426.00002
$ws.Range("M113").Value = -51723020
$ws.Range("N113").Value = -5766.000019999999
$ws.Range("H131").Value = 858.14
$ws.Range("J131").Value = 912.75
$ws.Range("L131").Value = 2738.25
$ws.Range("N131").Value = -12818.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 6666.6665
$ws.Range("J49").Value = 6666.6665
$ws.Range("L49").Value = 6666.6665
$ws.Range("N49").Value = -7034.6665
$ws.Range("H74").Value = 19998
$ws.Range("J74").Value = 19998
$ws.Range("L74").Value = 19998
$ws.Range("N74").Value = -21870
$ws.Range("H77").Value = 19998
$ws.Range("J77").Value = 19998
$ws.Range("L77").Value = 59994
$ws.Range("N77").Value = -69354
$ws.Range("H97").Value = 806
$ws.Range("I97").Value = 678.5
$ws.Range("J97").Value = 1112
$ws.Range("K97").Value = 678.5
$ws.Range("L97").Value = 1112
$ws.Range("M97").Value = -182.5
$ws.Range("N97").Value = -2104

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2403.238
$ws.Range("I7").Value = 2403.238
$ws.Range("K7").Value = 2403.238
$ws.Range("M7").Value = -2291.238
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 5000
$ws.Range("N42").Value = -6126
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5294
$ws.Range("H68").Value = 2317.0833
$ws.Range("I68").Value = 2700.4
$ws.Range("J68").Value = 2043.2858
$ws.Range("K68").Value = 2700.4
$ws.Range("L68").Value = 2043.2858
$ws.Range("M68").Value = -1951.4
$ws.Range("N68").Value = -3541.2858
$ws.Range("H71").Value = 2317.0833
$ws.Range("I71").Value = 2700.4
$ws.Range("J71").Value = 2043.2858
$ws.Range("K71").Value = 13502
$ws.Range("L71").Value = 10216.429
$ws.Range("M71").Value = -9758
$ws.Range("N71").Value = -17704.429
$ws.Range("H93").Value = 13856.75
$ws.Range("I93").Value = 26312.5
$ws.Range("J93").Value = 1401
$ws.Range("K93").Value = 26312.5
$ws.Range("L93").Value = 1401
$ws.Range("M93").Value = -25064.5
$ws.Range("N93").Value = -3897
$ws.Range("H100").Value = 85002750
$ws.Range("I100").Value = 2503375
$ws.Range("J100").Value = 250001500
$ws.Range("K100").Value = 2503375
$ws.Range("L100").Value = 250001500
$ws.Range("M100").Value = -2502834
$ws.Range("N100").Value = -250002582
$ws.Range("H126").Value = 2403.238
$ws.Range("I126").Value = 2403.238
$ws.Range("K126").Value = 7209.714
$ws.Range("M126").Value = -4739.714
$ws.Range("H132").Value = 1654.83
$ws.Range("I132").Value = 1647.7858
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4943.357400000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2413.357400000001
$ws.Range("N132").Value = -11060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2018.5
$ws.Range("J96").Value = 2018.5
$ws.Range("L96").Value = 2018.5
$ws.Range("N96").Value = -4764.5
$ws.Range("H100").Value = 562.44446
$ws.Range("I100").Value = 508.85715
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 1017.7143
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -476.7143
$ws.Range("N100").Value = -2582
$ws.Range("H132").Value = 1750.5625
$ws.Range("I132").Value = 1031.8148
$ws.Range("J132").Value = 2674.6667
$ws.Range("K132").Value = 3095.4444
$ws.Range("L132").Value = 8024.000100000001
$ws.Range("M132").Value = -565.4444000000003
$ws.Range("N132").Value = -13084.0001
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360
